$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking score for correct answer (row 11 "Marking")
$ws.Range("B11").Value = 5

# Update total correct marks (row 12 "Total")
$ws.Range("B12").Value = 85

# Update correct/total marks summary text
$ws.Range("E12").Value = "85/140"
